$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$new1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$new1.Name = "IbgCCSameDayInputter"

# temp sheets to consume sheetIds
$tmp = $wb.Worksheets.Add()
$tmp2 = $wb.Worksheets.Add()
$tmp.Delete()
$tmp2.Delete()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$new2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$new2.Name = "IbgCCNormalDayInputter"

Write-Host $wb.Worksheets.Count
